$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Header cell G1 - copy formatting (style) from F1, then set the text
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "language"

# Data cells G2:G26 - programming language per user
$ws.Range("G2").Value = "C++"
$ws.Range("G3").Value = "C++"
$ws.Range("G4").Value = "C++"
$ws.Range("G5").Value = "Python"
$ws.Range("G6").Value = "Python3"
$ws.Range("G7").Value = "Python"
$ws.Range("G8").Value = "C++"
$ws.Range("G9").Value = "C++"
$ws.Range("G10").Value = "C++"
$ws.Range("G11").Value = "C++"
$ws.Range("G12").Value = "C++"
$ws.Range("G13").Value = "C++"
$ws.Range("G14").Value = "Python3"
$ws.Range("G15").Value = "C++"
$ws.Range("G16").Value = "C++"
$ws.Range("G17").Value = "JavaScript"
$ws.Range("G18").Value = "Python3"
$ws.Range("G19").Value = "C++"
$ws.Range("G20").Value = "Java"
$ws.Range("G21").Value = "C++"
$ws.Range("G22").Value = "C++"
$ws.Range("G23").Value = "Python"
$ws.Range("G24").Value = "Python3"
$ws.Range("G25").Value = "C++"
$ws.Range("G26").Value = "Python3"
